$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw result as row 87, keeping every value a plain text
# string (matching the existing rows, which all use t="str"). A leading
# apostrophe forces Excel to treat the date-like / numeric-like strings
# ("2025-12-12" and "251212") as text instead of inferring a date serial
# or a number.
$ws.Range("A87").Value = "'2025-12-12"
$ws.Range("B87").Value = "Pick 4"
$ws.Range("C87").Value = "'251212"
$ws.Range("D87").Value = "9-7-0-2"
$ws.Range("E87").Value = "2025-12-12T21:41:58.649+04:00"

# The apostrophe trick stamps a "quote prefix" style on A87/C87. Copy the
# formatting (only) from the row above so the new row keeps the same
# (default) cell style as the rest of the table.
$ws.Range("A86:E86").Copy()
$ws.Range("A87").PasteSpecial(-4122)
